$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("tasas")

$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.82 = 6630.96 pesos`n✅ 6630.96 pesos = 1.81 = 946.06 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws2.Range("N10").Value = 550
$ws2.Range("O10").Value = 3647.03
$ws2.Range("N12").Value = 3661.23
$ws2.Range("O12").Value = 522.3579999999999
